# "save data done + era data updated"
# Adds a new "Save" column (H) to the sheet: a bold/bordered header in H1
# (matching the existing header formatting used for B1:G1) plus the
# save-flag values for each data row (2-6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1: copy formatting from the neighboring "sum" header (G1)
# so it picks up the same bold font / border / alignment style, then set
# its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New "Save" data column values for rows 2-6
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 1
